# Insert two new rows (cycle / cycleStart) right after the "title" row,
# pushing all subsequent key/value rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 4 (before the current "humidity" row).
$ws.Rows("4:5").Insert()

# Populate the newly inserted rows (order matches the shared-strings table
# layout: cycle, cycleStart, BEGIN, CYCLE {0}).
$ws.Range("A4").Value = "cycle"
$ws.Range("A5").Value = "cycleStart"
$ws.Range("B5").Value = "BEGIN"
$ws.Range("B4").Value = "CYCLE {0}"

# Update the selection to match the target state.
$ws.Range("B5").Select()
